$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "100"
$ws.Range("D1").Value = 5
$ws.Range("Q19").Select()
